# 1. Update the "总计" (summary) sheet: insert new row for 2022-Q3 at the top of the data block
$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# Capture the existing quarter rows (rows 2..8) before they get shifted down
$existingRows = @()
for ($r = 2; $r -le 8; $r++) {
    $label = $summary.Cells.Item($r,2).Value()
    $count = $summary.Cells.Item($r,3).Value()
    $value = $summary.Cells.Item($r,4).Value()
    $existingRows += ,@($label,$count,$value)
}

# Extend the index-column formatting (style) down to the new last row (row 9)
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)

# Re-write the old rows shifted down by one row (now rows 3..9)
for ($i = 0; $i -lt $existingRows.Count; $i++) {
    $r = $i + 3
    $summary.Cells.Item($r,1).Value = $i
    $summary.Cells.Item($r,2).Value = $existingRows[$i][0]
    $summary.Cells.Item($r,3).Value = $existingRows[$i][1]
    $summary.Cells.Item($r,4).Value = $existingRows[$i][2]
}

# Fill in the brand-new 2022-Q3 row at row 2
$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 20
$summary.Cells.Item(2,4).Value = 2.82

# 2. Insert a brand-new "2022-Q3" worksheet right after "总计" and before "2022-Q2"
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($q2Sheet, [System.Reflection.Missing]::Value)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"
$q3Sheet.Cells.ClearContents()

# Header row
$q3Sheet.Cells.Item(1,2).Value = "基金代码"
$q3Sheet.Cells.Item(1,3).Value = "基金名称"
$q3Sheet.Cells.Item(1,4).Value = "基金规模"
$q3Sheet.Cells.Item(1,5).Value = "股票总仓位"
$q3Sheet.Cells.Item(1,6).Value = "仓位占比"
$q3Sheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3Sheet.Cells.Item(1,8).Value = "仓位排名"

# Data rows
$q3Sheet.Cells.Item(2,1).Value = 0
$q3Sheet.Cells.Item(2,2).NumberFormat = "@"
$q3Sheet.Cells.Item(2,2).Value = "960010"
$q3Sheet.Cells.Item(2,3).Value = "工银核心价值混合H"
$q3Sheet.Cells.Item(2,4).NumberFormat = "@"
$q3Sheet.Cells.Item(2,4).Value = "49.15"
$q3Sheet.Cells.Item(2,5).NumberFormat = "@"
$q3Sheet.Cells.Item(2,5).Value = "85.16"
$q3Sheet.Cells.Item(2,6).NumberFormat = "@"
$q3Sheet.Cells.Item(2,6).Value = "2.87"
$q3Sheet.Cells.Item(2,7).NumberFormat = "@"
$q3Sheet.Cells.Item(2,7).Value = "1.4106"
$q3Sheet.Cells.Item(2,8).Value = 8

$q3Sheet.Cells.Item(3,1).Value = 1
$q3Sheet.Cells.Item(3,2).NumberFormat = "@"
$q3Sheet.Cells.Item(3,2).Value = "001719"
$q3Sheet.Cells.Item(3,3).Value = "工银国家战略股票"
$q3Sheet.Cells.Item(3,4).NumberFormat = "@"
$q3Sheet.Cells.Item(3,4).Value = "6.95"
$q3Sheet.Cells.Item(3,5).NumberFormat = "@"
$q3Sheet.Cells.Item(3,5).Value = "93.98"
$q3Sheet.Cells.Item(3,6).NumberFormat = "@"
$q3Sheet.Cells.Item(3,6).Value = "5.03"
$q3Sheet.Cells.Item(3,7).NumberFormat = "@"
$q3Sheet.Cells.Item(3,7).Value = "0.3496"
$q3Sheet.Cells.Item(3,8).Value = 7

$q3Sheet.Cells.Item(4,1).Value = 2
$q3Sheet.Cells.Item(4,2).NumberFormat = "@"
$q3Sheet.Cells.Item(4,2).Value = "159745"
$q3Sheet.Cells.Item(4,3).Value = "国泰中证全指建筑材料ETF"
$q3Sheet.Cells.Item(4,4).NumberFormat = "@"
$q3Sheet.Cells.Item(4,4).Value = "7.92"
$q3Sheet.Cells.Item(4,5).NumberFormat = "@"
$q3Sheet.Cells.Item(4,5).Value = "99.14"
$q3Sheet.Cells.Item(4,6).NumberFormat = "@"
$q3Sheet.Cells.Item(4,6).Value = "3.64"
$q3Sheet.Cells.Item(4,7).NumberFormat = "@"
$q3Sheet.Cells.Item(4,7).Value = "0.2883"
$q3Sheet.Cells.Item(4,8).Value = 6

$q3Sheet.Cells.Item(5,1).Value = 3
$q3Sheet.Cells.Item(5,2).NumberFormat = "@"
$q3Sheet.Cells.Item(5,2).Value = "001008"
$q3Sheet.Cells.Item(5,3).Value = "工银国企改革主题股票"
$q3Sheet.Cells.Item(5,4).NumberFormat = "@"
$q3Sheet.Cells.Item(5,4).Value = "8.66"
$q3Sheet.Cells.Item(5,5).NumberFormat = "@"
$q3Sheet.Cells.Item(5,5).Value = "81.53"
$q3Sheet.Cells.Item(5,6).NumberFormat = "@"
$q3Sheet.Cells.Item(5,6).Value = "3.21"
$q3Sheet.Cells.Item(5,7).NumberFormat = "@"
$q3Sheet.Cells.Item(5,7).Value = "0.2780"
$q3Sheet.Cells.Item(5,8).Value = 6

$q3Sheet.Cells.Item(6,1).Value = 4
$q3Sheet.Cells.Item(6,2).NumberFormat = "@"
$q3Sheet.Cells.Item(6,2).Value = "004856"
$q3Sheet.Cells.Item(6,3).Value = "广发中证全指建筑材料指数A"
$q3Sheet.Cells.Item(6,4).NumberFormat = "@"
$q3Sheet.Cells.Item(6,4).Value = "7.66"
$q3Sheet.Cells.Item(6,5).NumberFormat = "@"
$q3Sheet.Cells.Item(6,5).Value = "93.74"
$q3Sheet.Cells.Item(6,6).NumberFormat = "@"
$q3Sheet.Cells.Item(6,6).Value = "3.44"
$q3Sheet.Cells.Item(6,7).NumberFormat = "@"
$q3Sheet.Cells.Item(6,7).Value = "0.2635"
$q3Sheet.Cells.Item(6,8).Value = 6

$q3Sheet.Cells.Item(7,1).Value = 5
$q3Sheet.Cells.Item(7,2).NumberFormat = "@"
$q3Sheet.Cells.Item(7,2).Value = "004857"
$q3Sheet.Cells.Item(7,3).Value = "广发中证全指建筑材料指数C"
$q3Sheet.Cells.Item(7,4).NumberFormat = "@"
$q3Sheet.Cells.Item(7,4).Value = "6.12"
$q3Sheet.Cells.Item(7,5).NumberFormat = "@"
$q3Sheet.Cells.Item(7,5).Value = "93.74"
$q3Sheet.Cells.Item(7,6).NumberFormat = "@"
$q3Sheet.Cells.Item(7,6).Value = "3.44"
$q3Sheet.Cells.Item(7,7).NumberFormat = "@"
$q3Sheet.Cells.Item(7,7).Value = "0.2105"
$q3Sheet.Cells.Item(7,8).Value = 6

$q3Sheet.Cells.Item(8,1).Value = 6
$q3Sheet.Cells.Item(8,2).NumberFormat = "@"
$q3Sheet.Cells.Item(8,2).Value = "004405"
$q3Sheet.Cells.Item(8,3).Value = "国寿安保稳寿混合A"
$q3Sheet.Cells.Item(8,4).NumberFormat = "@"
$q3Sheet.Cells.Item(8,4).Value = "4.69"
$q3Sheet.Cells.Item(8,5).NumberFormat = "@"
$q3Sheet.Cells.Item(8,5).Value = "24.56"
$q3Sheet.Cells.Item(8,6).NumberFormat = "@"
$q3Sheet.Cells.Item(8,6).Value = "1.10"
$q3Sheet.Cells.Item(8,7).NumberFormat = "@"
$q3Sheet.Cells.Item(8,7).Value = "0.0516"
$q3Sheet.Cells.Item(8,8).Value = 2

$q3Sheet.Cells.Item(9,1).Value = 7
$q3Sheet.Cells.Item(9,2).NumberFormat = "@"
$q3Sheet.Cells.Item(9,2).Value = "516750"
$q3Sheet.Cells.Item(9,3).Value = "富国中证全指建筑材料ETF"
$q3Sheet.Cells.Item(9,4).NumberFormat = "@"
$q3Sheet.Cells.Item(9,4).Value = "0.82"
$q3Sheet.Cells.Item(9,5).NumberFormat = "@"
$q3Sheet.Cells.Item(9,5).Value = "98.46"
$q3Sheet.Cells.Item(9,6).NumberFormat = "@"
$q3Sheet.Cells.Item(9,6).Value = "3.59"
$q3Sheet.Cells.Item(9,7).NumberFormat = "@"
$q3Sheet.Cells.Item(9,7).Value = "0.0294"
$q3Sheet.Cells.Item(9,8).Value = 7

$q3Sheet.Cells.Item(10,1).Value = 8
$q3Sheet.Cells.Item(10,2).NumberFormat = "@"
$q3Sheet.Cells.Item(10,2).Value = "011073"
$q3Sheet.Cells.Item(10,3).Value = "鹏华安润混合A"
$q3Sheet.Cells.Item(10,4).NumberFormat = "@"
$q3Sheet.Cells.Item(10,4).Value = "1.48"
$q3Sheet.Cells.Item(10,5).NumberFormat = "@"
$q3Sheet.Cells.Item(10,5).Value = "25.50"
$q3Sheet.Cells.Item(10,6).NumberFormat = "@"
$q3Sheet.Cells.Item(10,6).Value = "1.63"
$q3Sheet.Cells.Item(10,7).NumberFormat = "@"
$q3Sheet.Cells.Item(10,7).Value = "0.0241"
$q3Sheet.Cells.Item(10,8).Value = 4

$q3Sheet.Cells.Item(11,1).Value = 9
$q3Sheet.Cells.Item(11,2).NumberFormat = "@"
$q3Sheet.Cells.Item(11,2).Value = "012419"
$q3Sheet.Cells.Item(11,3).Value = "天弘国证建材指数C"
$q3Sheet.Cells.Item(11,4).NumberFormat = "@"
$q3Sheet.Cells.Item(11,4).Value = "0.63"
$q3Sheet.Cells.Item(11,5).NumberFormat = "@"
$q3Sheet.Cells.Item(11,5).Value = "94.93"
$q3Sheet.Cells.Item(11,6).NumberFormat = "@"
$q3Sheet.Cells.Item(11,6).Value = "3.09"
$q3Sheet.Cells.Item(11,7).NumberFormat = "@"
$q3Sheet.Cells.Item(11,7).Value = "0.0195"
$q3Sheet.Cells.Item(11,8).Value = 9

$q3Sheet.Cells.Item(12,1).Value = 10
$q3Sheet.Cells.Item(12,2).NumberFormat = "@"
$q3Sheet.Cells.Item(12,2).Value = "011761"
$q3Sheet.Cells.Item(12,3).Value = "平安鑫瑞混合A"
$q3Sheet.Cells.Item(12,4).NumberFormat = "@"
$q3Sheet.Cells.Item(12,4).Value = "0.59"
$q3Sheet.Cells.Item(12,5).NumberFormat = "@"
$q3Sheet.Cells.Item(12,5).Value = "23.80"
$q3Sheet.Cells.Item(12,6).NumberFormat = "@"
$q3Sheet.Cells.Item(12,6).Value = "1.56"
$q3Sheet.Cells.Item(12,7).NumberFormat = "@"
$q3Sheet.Cells.Item(12,7).Value = "0.0092"
$q3Sheet.Cells.Item(12,8).Value = 2

$q3Sheet.Cells.Item(13,1).Value = 11
$q3Sheet.Cells.Item(13,2).NumberFormat = "@"
$q3Sheet.Cells.Item(13,2).Value = "512590"
$q3Sheet.Cells.Item(13,3).Value = "浦银安盛中证高股息精选ETF"
$q3Sheet.Cells.Item(13,4).NumberFormat = "@"
$q3Sheet.Cells.Item(13,4).Value = "0.45"
$q3Sheet.Cells.Item(13,5).NumberFormat = "@"
$q3Sheet.Cells.Item(13,5).Value = "90.87"
$q3Sheet.Cells.Item(13,6).NumberFormat = "@"
$q3Sheet.Cells.Item(13,6).Value = "2.01"
$q3Sheet.Cells.Item(13,7).NumberFormat = "@"
$q3Sheet.Cells.Item(13,7).Value = "0.0090"
$q3Sheet.Cells.Item(13,8).Value = 6

$q3Sheet.Cells.Item(14,1).Value = 12
$q3Sheet.Cells.Item(14,2).NumberFormat = "@"
$q3Sheet.Cells.Item(14,2).Value = "004406"
$q3Sheet.Cells.Item(14,3).Value = "国寿安保稳寿混合C"
$q3Sheet.Cells.Item(14,4).NumberFormat = "@"
$q3Sheet.Cells.Item(14,4).Value = "0.62"
$q3Sheet.Cells.Item(14,5).NumberFormat = "@"
$q3Sheet.Cells.Item(14,5).Value = "24.56"
$q3Sheet.Cells.Item(14,6).NumberFormat = "@"
$q3Sheet.Cells.Item(14,6).Value = "1.10"
$q3Sheet.Cells.Item(14,7).NumberFormat = "@"
$q3Sheet.Cells.Item(14,7).Value = "0.0068"
$q3Sheet.Cells.Item(14,8).Value = 2

$q3Sheet.Cells.Item(15,1).Value = 13
$q3Sheet.Cells.Item(15,2).NumberFormat = "@"
$q3Sheet.Cells.Item(15,2).Value = "159787"
$q3Sheet.Cells.Item(15,3).Value = "易方达中证全指建筑材料ETF"
$q3Sheet.Cells.Item(15,4).NumberFormat = "@"
$q3Sheet.Cells.Item(15,4).Value = "0.17"
$q3Sheet.Cells.Item(15,5).NumberFormat = "@"
$q3Sheet.Cells.Item(15,5).Value = "94.24"
$q3Sheet.Cells.Item(15,6).NumberFormat = "@"
$q3Sheet.Cells.Item(15,6).Value = "3.56"
$q3Sheet.Cells.Item(15,7).NumberFormat = "@"
$q3Sheet.Cells.Item(15,7).Value = "0.0061"
$q3Sheet.Cells.Item(15,8).Value = 6

$q3Sheet.Cells.Item(16,1).Value = 14
$q3Sheet.Cells.Item(16,2).NumberFormat = "@"
$q3Sheet.Cells.Item(16,2).Value = "011762"
$q3Sheet.Cells.Item(16,3).Value = "平安鑫瑞混合C"
$q3Sheet.Cells.Item(16,4).NumberFormat = "@"
$q3Sheet.Cells.Item(16,4).Value = "0.28"
$q3Sheet.Cells.Item(16,5).NumberFormat = "@"
$q3Sheet.Cells.Item(16,5).Value = "23.80"
$q3Sheet.Cells.Item(16,6).NumberFormat = "@"
$q3Sheet.Cells.Item(16,6).Value = "1.56"
$q3Sheet.Cells.Item(16,7).NumberFormat = "@"
$q3Sheet.Cells.Item(16,7).Value = "0.0044"
$q3Sheet.Cells.Item(16,8).Value = 2

$q3Sheet.Cells.Item(17,1).Value = 15
$q3Sheet.Cells.Item(17,2).NumberFormat = "@"
$q3Sheet.Cells.Item(17,2).Value = "012405"
$q3Sheet.Cells.Item(17,3).Value = "天弘国证建材指数A"
$q3Sheet.Cells.Item(17,4).NumberFormat = "@"
$q3Sheet.Cells.Item(17,4).Value = "0.13"
$q3Sheet.Cells.Item(17,5).NumberFormat = "@"
$q3Sheet.Cells.Item(17,5).Value = "94.93"
$q3Sheet.Cells.Item(17,6).NumberFormat = "@"
$q3Sheet.Cells.Item(17,6).Value = "3.09"
$q3Sheet.Cells.Item(17,7).NumberFormat = "@"
$q3Sheet.Cells.Item(17,7).Value = "0.0040"
$q3Sheet.Cells.Item(17,8).Value = 9

$q3Sheet.Cells.Item(18,1).Value = 16
$q3Sheet.Cells.Item(18,2).NumberFormat = "@"
$q3Sheet.Cells.Item(18,2).Value = "011074"
$q3Sheet.Cells.Item(18,3).Value = "鹏华安润混合C"
$q3Sheet.Cells.Item(18,4).NumberFormat = "@"
$q3Sheet.Cells.Item(18,4).Value = "0.19"
$q3Sheet.Cells.Item(18,5).NumberFormat = "@"
$q3Sheet.Cells.Item(18,5).Value = "25.50"
$q3Sheet.Cells.Item(18,6).NumberFormat = "@"
$q3Sheet.Cells.Item(18,6).Value = "1.63"
$q3Sheet.Cells.Item(18,7).NumberFormat = "@"
$q3Sheet.Cells.Item(18,7).Value = "0.0031"
$q3Sheet.Cells.Item(18,8).Value = 4

$q3Sheet.Cells.Item(19,1).Value = 17
$q3Sheet.Cells.Item(19,2).NumberFormat = "@"
$q3Sheet.Cells.Item(19,2).Value = "008719"
$q3Sheet.Cells.Item(19,3).Value = "德邦安顺混合A"
$q3Sheet.Cells.Item(19,4).NumberFormat = "@"
$q3Sheet.Cells.Item(19,4).Value = "0.29"
$q3Sheet.Cells.Item(19,5).NumberFormat = "@"
$q3Sheet.Cells.Item(19,5).Value = "22.84"
$q3Sheet.Cells.Item(19,6).NumberFormat = "@"
$q3Sheet.Cells.Item(19,6).Value = "1.00"
$q3Sheet.Cells.Item(19,7).NumberFormat = "@"
$q3Sheet.Cells.Item(19,7).Value = "0.0029"
$q3Sheet.Cells.Item(19,8).Value = 5

$q3Sheet.Cells.Item(20,1).Value = 18
$q3Sheet.Cells.Item(20,2).NumberFormat = "@"
$q3Sheet.Cells.Item(20,2).Value = "008720"
$q3Sheet.Cells.Item(20,3).Value = "德邦安顺混合C"
$q3Sheet.Cells.Item(20,4).NumberFormat = "@"
$q3Sheet.Cells.Item(20,4).Value = "0.26"
$q3Sheet.Cells.Item(20,5).NumberFormat = "@"
$q3Sheet.Cells.Item(20,5).Value = "22.84"
$q3Sheet.Cells.Item(20,6).NumberFormat = "@"
$q3Sheet.Cells.Item(20,6).Value = "1.00"
$q3Sheet.Cells.Item(20,7).NumberFormat = "@"
$q3Sheet.Cells.Item(20,7).Value = "0.0026"
$q3Sheet.Cells.Item(20,8).Value = 5

$q3Sheet.Cells.Item(21,1).Value = 19
$q3Sheet.Cells.Item(21,2).NumberFormat = "@"
$q3Sheet.Cells.Item(21,2).Value = "481001"
$q3Sheet.Cells.Item(21,3).Value = "工银核心价值混合A"
$q3Sheet.Cells.Item(21,4).NumberFormat = "@"
$q3Sheet.Cells.Item(21,4).Value = "-5.40"
$q3Sheet.Cells.Item(21,5).NumberFormat = "@"
$q3Sheet.Cells.Item(21,5).Value = "85.16"
$q3Sheet.Cells.Item(21,6).NumberFormat = "@"
$q3Sheet.Cells.Item(21,6).Value = "2.87"
$q3Sheet.Cells.Item(21,7).NumberFormat = "@"
$q3Sheet.Cells.Item(21,7).Value = "-0.1550"
$q3Sheet.Cells.Item(21,8).Value = 8

# Extend the index-column (A) formatting down past the templates original last row (19) to the new last row (21)
$q3Sheet.Range("A19").Copy()
$q3Sheet.Range("A20:A21").PasteSpecial(-4122)

Write-Host "Added 2022-Q3 sheet and updated summary sheet"
